$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two rows that are no longer present (B10-0943 at row 11 and
# B12-0557 at row 12). Deleting them shifts everything below up by two rows,
# matching the new layout where B12-0393-2 becomes row 11.
$ws.Rows("11:12").Delete()

# Append the new rows at the bottom of the list (rows 24-29 in the final
# layout: one text id, one numeric id, and four more text ids).
$ws.Range("A24").Value = "2308_GCA_000054005-2"
$ws.Range("A25").Value = 2002734562
$ws.Range("A26").Value = "B17-0712"
$ws.Range("A27").Value = "B17-0712-RETRO"
$ws.Range("A28").Value = "B17-0712-SMAM"
$ws.Range("A29").Value = "B17-0712-TRACH"

# Update the active selection to match the post-edit state (cursor parked on
# the first empty row below the data).
[void]$ws.Range("A30").Select()
